$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.003.83"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.557.28"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'207.38"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "'21.67"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.559.45"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "27.009.77"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "'216.23"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("D24").Value = "'1.96"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").Value = "'152.61"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "'6.66"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").Value = "'0.0463"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("D33").Value = "1.402.42"
$ws.Range("E33").Value = "  +4.15%  "
$ws.Range("E34").Value = "  +3.42%  "
$ws.Range("E35").Value = "  +3.82%  "
$ws.Range("D36").Value = "'0.960"
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("D39").Value = "'0.524"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").Value = "'0.991"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'2.28"
$ws.Range("E43").Value = "  +4.03%  "
$ws.Range("D44").Value = "'5.47"
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("D45").Value = "'63.98"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "1.693.82"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "'0.0959"
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("E51").Value = "  +0.44%  "
